# Update "实验结果记录" workbook: add a new experiment result row (row 16)
# for the "SD-SHW" (step-distance reward + summed hit-wall penalty) run,
# and refresh the active window's zoom/selection to match the author's
# last-saved view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# ---- New row 16: short name / detail / per-map observations ----------
$ws.Range("A16").Value = "SD-SHW"
$ws.Range("B16").Value = "PPO use step distance reward + multiply critic lr + train every episode + summed hit wall penalty vs. Random"
$ws.Range("C16").Value = "走的不直线，所以有时输"
$ws.Range("D16").Value = "撞墙还是不少"
$ws.Range("E16").Value = "撞墙不多，但是有干扰的话就不知所措了"
$ws.Range("J16").Value = "很强，马上就过"
$ws.Range("K16").Value = "受干扰的话可能寄"

# ---- Match formatting used by the rest of the results table -----------
# Column A uses the centered "short name" style.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)

# Column B (details) carries no special fill/alignment.
$ws.Range("B15").Copy()
$ws.Range("B16").PasteSpecial(-4122)

# Result columns (map observations) use the light-green highlight style.
$ws.Range("C15").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C15").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("C15").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("C15").Copy()
$ws.Range("J16").PasteSpecial(-4122)
$ws.Range("C15").Copy()
$ws.Range("K16").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---- Window view: restore zoom level & active cell selection ----------
$excel.ActiveWindow.Zoom = 118
$ws.Range("I12").Select()
